$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking") updates
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total") updates
$ws.Range("B12").Value = 80
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "76 / 112"
